# تعديل تلقائي في شيت Card4 by admin at 2025-12-06 18:33:52
# Update the "card" identifier in column A from 2 to 4 for the Card4 sheet
# (rows 3-7 and 9-13; row 8 already reads 4).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card4")

$rows = @(3, 4, 5, 6, 7, 9, 10, 11, 12, 13)
foreach ($r in $rows) {
    # Leading apostrophe keeps the value stored as text ("4"), matching the
    # existing text-typed "card" values in column A instead of coercing to
    # a number.
    $ws.Cells.Item($r, 1).Value = "'4"
}
